$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New TPM-based data replaces the previous 9-row result with 6 rows
# (the ECs sending-cluster rows are no longer present).
$data = @(
    @("FAPs", "Ctf1", "Il6st", "ECs", 2, 0.6666666666666666, 0.7996153333333332, 2.398846, 0.373961854065697, 0.3739618540656971, 3, 1, 20.854426, 62.563278, 0.1507164072139519, 0.1507164072139519, 16.67551879746533, 150.079669177188, 0.05636218707985007, 0.05636218707985007),
    @("FAPs", "Ctf1", "Il6st", "FAPs", 2, 0.6666666666666666, 0.7996153333333332, 2.398846, 0.373961854065697, 0.3739618540656971, 3, 1, 94.96115633333334, 284.883469, 0.6862909728343718, 0.6862909728343718, 75.9323966751971, 683.391570076774, 0.2566466446296926, 0.2566466446296926),
    @("FAPs", "Ctf1", "Il6st", "MuSCs", 2, 0.6666666666666666, 0.7996153333333332, 2.398846, 0.373961854065697, 0.3739618540656971, 3, 1, 22.553069, 67.65920700000001, 0.1629926199516763, 0.1629926199516763, 18.03377978612467, 162.304018075122, 0.06095302235615438, 0.06095302235615439),
    @("MuSCs", "Ctf1", "Il6st", "ECs", 3, 1, 1.338611666666667, 4.015835, 0.6260381459343028, 0.626038145934303, 3, 1, 20.854426, 62.563278, 0.1507164072139519, 0.1507164072139519, 27.91597794523667, 251.24380150713, 0.09435422013410186, 0.09435422013410187),
    @("MuSCs", "Ctf1", "Il6st", "FAPs", 3, 1, 1.338611666666667, 4.015835, 0.6260381459343028, 0.626038145934303, 3, 1, 94.96115633333334, 284.883469, 0.6862909728343718, 0.6862909728343718, 127.1161117479572, 1144.045005731615, 0.4296443282046791, 0.4296443282046792),
    @("MuSCs", "Ctf1", "Il6st", "MuSCs", 3, 1, 1.338611666666667, 4.015835, 0.6260381459343028, 0.626038145934303, 3, 1, 22.553069, 67.65920700000001, 0.1629926199516763, 0.1629926199516763, 30.18980128253834, 271.708211542845, 0.1020395975955219, 0.1020395975955219)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}

# Drop the now-unused trailing rows (previously rows 8-10).
$ws.Range("A8:T10").ClearContents()

